# Rewrites the final bullet list ("Allgemeine DOTO") with the new set of
# TODO items (with priorities), per commit "Neue TODOs (mit prio)".
#
# The old list (numId=18, all ilvl=0) was:
#   Javascript Array usw.
#   Prototyping
#   Onblur und Co, Funktionen
#   WebWroker
#   Callback
#   Validation (aber wichtig für Arbeit)
#   CSS insbesondere der CSS3 advanced Stuff
#
# It is replaced by a 21-item list (mixing ilvl 0 / 1, and a couple of
# sub-items highlighted in red).

$d = $word.ActiveDocument

# Locate the first paragraph of the old list by scanning (stable even if
# indices were to shift) rather than hard-coding its paragraph number.
$n = $d.Paragraphs.Count
$firstIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Javascript Array usw.") {
        $firstIndex = $i
        break
    }
}
if ($firstIndex -eq -1) {
    throw "Could not locate the start of the 'Allgemeine DOTO' list"
}

# The new content for each bullet. Each entry is a list of (text, color)
# run tuples -- color is $null for "automatic" (no explicit color).
$items = New-Object System.Collections.ArrayList

function Add-Item($ilvl, $runs) {
    [void]$items.Add(@{ ilvl = $ilvl; runs = $runs })
}

Add-Item 0 @(@{ t = "Tabbellen"; c = $null })
Add-Item 0 @(@{ t = "Mit JavaScript (JQuery) Elemente verschieben"; c = $null })
Add-Item 0 @(@{ t = "Validation (aber wichtig für Arbeiten)"; c = $null })
Add-Item 0 @(@{ t = "`"Programmatically add and modify HTML elements; implement media controls; implement HTML5 canvas and SVG graphics"; c = $null })
Add-Item 1 @(@{ t = "DOM mitlerweile klar, auch JQuery"; c = $null })
Add-Item 1 @(@{ t = "HTML elemente hinzufügen"; c = 255 }, @{ t = " (mit JQuery), noch üben"; c = $null })
Add-Item 1 @(@{ t = "HTML5 canvas uns SVG offen"; c = 255 })
Add-Item 0 @(@{ t = "`"Implement storage APIs, AppCache API, and Geolocation API`""; c = $null })
Add-Item 1 @(@{ t = "storage verinnerlichen"; c = $null })
Add-Item 1 @(@{ t = "Rest Unterlagen prüfen"; c = $null })
Add-Item 0 @(@{ t = "Javascript Gültigkeitsbereiche (global, local, this)"; c = $null })
Add-Item 0 @(@{ t = "`"bubbled events`""; c = $null })
Add-Item 1 @(@{ t = "wann ist ein Event `"gebubbled`""; c = $null })
Add-Item 1 @(@{ t = "sind gefangene Events automatisch aus den `"bubble`" Vorgang entfernt"; c = $null })
Add-Item 0 @(@{ t = "das `"null`" Objekt (JS)"; c = $null })
Add-Item 0 @(@{ t = "AJAX Calls verinnerlichen"; c = $null })
Add-Item 0 @(@{ t = "Implementierung von WebWorkern (Notizen), Timeouts"; c = $null })
Add-Item 0 @(@{ t = "(De)Serialisierung von Daten (JSON, XML)"; c = $null })
Add-Item 0 @(@{ t = "Text-Styling"; c = $null })
Add-Item 0 @(@{ t = "CSS Box-Properties (runde Ecken ...)"; c = $null })
Add-Item 0 @(@{ t = "Prototyping"; c = $null })

$oldCount = 7
$newCount = $items.Count

# Grow the paragraph run so there is exactly one paragraph per new item.
# Each InsertParagraphAfter() on the (growing) last paragraph of the list
# clones its pPr (style + numPr), and becomes the new last paragraph.
$lastOldIndex = $firstIndex + $oldCount - 1
$toAdd = $newCount - $oldCount
for ($k = 0; $k -lt $toAdd; $k++) {
    $d.Paragraphs.Item($lastOldIndex + $k).Range.InsertParagraphAfter() | Out-Null
}

# Now paragraphs $firstIndex .. ($firstIndex + $newCount - 1) are the ones
# to fill in, one per item, all based on ListParagraph / numId 18.
for ($i = 0; $i -lt $newCount; $i++) {
    $item = $items[$i]
    $p = $d.Paragraphs.Item($firstIndex + $i)

    # Full plain text first (single run, no special formatting), then go
    # back and recolor the sub-ranges that need it. Doing it in this order
    # (type first, color after) keeps uncolored runs free of any explicit
    # <w:color> element.
    $fullText = ""
    foreach ($r in $item.runs) { $fullText += $r.t }
    $p.Range.Text = $fullText

    $p.Range.ListFormat.ListLevelNumber = $item.ilvl + 1

    $pos = $p.Range.Start
    foreach ($r in $item.runs) {
        $len = $r.t.Length
        if ($r.c -ne $null) {
            $sub = $d.Range($pos, $pos + $len)
            $sub.Font.Color = $r.c
        }
        $pos = $pos + $len
    }
}

Write-Output "Rewrote $newCount TODO bullet items (was $oldCount)."
